# Apply the "Fixed update to excel issue" change:
#  1. Rename the "Requested quantity" header on "Weekly Quantity" -> "Weekly_PO_Qty"
#  2. Rename the "Requested quantity" header on "Monthly Trend"  -> "Monthly_PO_Qty"
#  3. Add a new "PO Forecast" worksheet (after "Monthly Trend") containing the
#     ds / PO_Forecast / yhat_lower / yhat_upper forecast table.

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity header -------------------------------------------------
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend header ----------------------------------------------------
$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. New "PO Forecast" sheet, inserted after the last existing sheet --------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"
$wsForecast.Range("A1:D1").Font.Bold = $true
$wsForecast.Range("A1:D1").Borders.LineStyle = 1
$wsForecast.Range("A1:D1").HorizontalAlignment = -4108
$wsForecast.Range("A1:D1").VerticalAlignment = -4160

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(44948.99999999999, 0,   -35.26047516093856, 15.45301455630845),
    @(44955.99999999999, 6,   -22.72923523415547, 30.52380391236834),
    @(44962.99999999999, 23,  -3.536202927279285, 50.35559478768705),
    @(44969.99999999999, 41,  12.90113439236435,  68.74716289468445),
    @(44990.99999999999, 92,  62.96342195097731,  119.6922867829692),
    @(44997.99999999999, 109, 79.13334440058102,  136.3172686392036),
    @(45004.99999999999, 126, 99.85283048676625,  152.8848800256834),
    @(45011.99999999999, 144, 117.4607148150848,  167.4693161334861),
    @(45018.99999999999, 161, 135.1220091625447,  187.0771681012736),
    @(45025.99999999999, 178, 153.5081252980575,  203.8299682107177),
    @(45032.99999999999, 195, 169.1230806395492,  221.9749711770414),
    @(45039.99999999999, 212, 185.8497282321856,  239.1947672691947),
    @(45046.99999999999, 229, 200.899640142415,   254.6087962034067)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Range("A$r").Value = $row[0]
    $wsForecast.Range("B$r").Value = $row[1]
    $wsForecast.Range("C$r").Value = $row[2]
    $wsForecast.Range("D$r").Value = $row[3]
    $r = $r + 1
}

$wsForecast.Range("A2:A14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
